# Commit: "Added numbering in the list topic wise"
#
# Column B holds a per-topic serial number (1, 2, 3, ...) for each
# contiguous block of rows that share the same topic in column A
# (e.g. "Matrix", "String", "Searching & Sorting", ...). Before this
# edit, column B was blank for rows 44-481; this fills it in, resetting
# the counter back to 1 at the start of each topic block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (firstRow, lastRow) for every topic block in rows 44-481.
$blocks = @(
    @(44, 53),
    @(56, 98),
    @(101, 136),
    @(139, 174),
    @(177, 211),
    @(214, 235),
    @(238, 272),
    @(275, 293),
    @(296, 333),
    @(336, 353),
    @(356, 399),
    @(402, 407),
    @(410, 469),
    @(472, 481)
)

foreach ($block in $blocks) {
    $firstRow = $block[0]
    $lastRow = $block[1]
    $n = 1
    for ($row = $firstRow; $row -le $lastRow; $row++) {
        $ws.Cells.Item($row, 2).Value = $n
        $n++
    }
}

# Real Excel re-wrapped row 53 (last row of the first block, "Matrix")
# to a custom height of 17.25 when its B cell was edited.
$ws.Rows(53).RowHeight = 17.25

# Final view state recorded in the saved file: scrolled further down the
# sheet, with B472:B481 (the last topic block) selected.
$ws.Activate()
$ws.Range("A271").Select()
$ws.Range("B472:B481").Select()
